# Add a new literature-matrix entry as row 7 on "Month 1".
#
# Columns (A..F): Title/Author, Paper, Problem statement, Key findings,
# Methodology, Relevance of literature  (column G "Theme" / H "Impact of
# Work" are left blank, matching the source diff which only populates A-F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titleAuthor   = "Anam Sajid, Haider Abbas, Kashif saleem"
$paper         = "Cloud-assisted IoT based SCADA systems security: A review of the state of the art and future challenges"
$problem       = "The IoT-cloud based SCADA ICS faces vulnerabilities as it is integrated with underlying legacy systems. Some of these vulnerabilities are as a result of the weak (insecure) communication protocols - Modbus, IEC 61850, and CI systems"
$keyFindings   = "The survey paper explores the different kinds of vulnerabilities faced since the adoption of the networked and the IoT cloud based SCADA system approaches. The paper reviews the different attacks that have hit several organizations siting CIA for security risk analysis. "
$methodology   = "SCADA systems are susceptible to vulnerabilities during communication and at the CI layer as well since most Industrial applications use commercial cloud services. This survey equally presents case scenarios of different threats: eavesdropping, man-in-the-middle attacks, data corruption owing to lack of proper security controls"
$relevance     = "This paper is quite important as it addresses the differen threats. It reviews the chronological advancement in the design and implementation of SCADA systems. The paper discusses the vulnerabilities of communication layer (overview) and proposes the use of multiple layers of security as well as redundancy as a preventative measure. "

# Plain-text cells A7:E7
$ws.Cells.Item(7, 1).Value = $titleAuthor
$ws.Cells.Item(7, 2).Value = $paper
$ws.Cells.Item(7, 3).Value = $problem
$ws.Cells.Item(7, 4).Value = $keyFindings
$ws.Cells.Item(7, 5).Value = $methodology

# Rich-text cell F7: "...communication layer (" + bold("overview") + ") and proposes..."
$ws.Cells.Item(7, 6).Value = $relevance

# Match the existing body-row formatting (Century Gothic 9pt, wrap, top-aligned)
for ($col = 1; $col -le 6; $col++) {
    $cell = $ws.Cells.Item(7, $col)
    $cell.Font.Name = "Century Gothic"
    $cell.Font.Size = 9
    $cell.VerticalAlignment = -4160   # xlTop
    $cell.WrapText = $true
}

# Bold the word "overview" inside F7, keeping the rest at the normal weight
$relIdx  = $relevance.IndexOf("overview") + 1
$relLen  = 8
$boldRun = $ws.Cells.Item(7, 6).Characters($relIdx, $relLen)
$boldRun.Font.Name = "Century Gothic"
$boldRun.Font.Size = 9
$boldRun.Font.Bold = $true

$afterIdx = $relIdx + $relLen
$afterLen = $relevance.Length - ($afterIdx - 1)
$afterRun = $ws.Cells.Item(7, 6).Characters($afterIdx, $afterLen)
$afterRun.Font.Name = "Century Gothic"
$afterRun.Font.Size = 9

# Row height matches the other wrapped body rows
$ws.Rows.Item(7).RowHeight = 211.2

# Update the view to mirror the saved selection/scroll position
$ws.Range("A7").Select()
